$wb = $excel.ActiveWorkbook

# --- Basic Options: add the "Template Version" row below the existing rows ---
$basic = $wb.Worksheets.Item("Basic Options")
$basic.Range("A7").Value = "Template Version"
$basic.Range("B7").Value = 2
$basic.Range("A8").Select() | Out-Null

# --- Add the new "Quantification" sheet after the last existing sheet ("Plots") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Quantification"

# Header row content (reuses existing shared strings "#Sample Group" / "Region")
$newSheet.Range("A1").Value = "#Sample Group"
$newSheet.Range("B1").Value = "Region"

# Match header formatting used on the other sheets (bold header row)
$newSheet.Range("A1:B1").Font.Bold = $true

# Autofit the first column, like the other sheets' headers
$newSheet.Columns.Item(1).AutoFit()

# Select the header row (mirrors the original author clicking the row header)
$newSheet.Rows.Item(1).Select() | Out-Null
